$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C11 value
$ws.Range("C11").Value = 0.3298611111111111

# Add new row 12 data
$ws.Range("A12").Value = 43340
$ws.Range("B12").Value = 0.55972222222222223
$ws.Range("C12").Value = 0.75347222222222221

# Add new row 13 data (date only)
$ws.Range("A13").Value = 43341

# Update selection to match diff
$ws.Range("B13").Select()
